$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Duplicate the "2021-Q4" sheet (same column layout as the new
#    quarter) and place the copy right before "总计", then rename it
#    to "2022-Q1".
# ------------------------------------------------------------------
$sheetQ4    = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetQ4.Copy($sheetTotal)

# Re-fetch the freshly created sheet by its (stable) position instead
# of by name - worksheet references can go stale after the sheet
# collection is mutated.
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "2022-Q1"

# ------------------------------------------------------------------
# 2) Update the fund data on the new "2022-Q1" sheet.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")

# D2:G2 hold numeric-looking figures that must stay TEXT (as in the
# other quarter sheets), so force a text format before assigning the
# values, then restore the default "Normal" style so no stray number
# format lingers on the cells.
$q1.Range("D2:G2").NumberFormat = "@"
$q1.Range("D2").Value = "297.64"
$q1.Range("E2").Value = "57.54"
$q1.Range("F2").Value = "1.99"
$q1.Range("G2").Value = "5.9230"
$q1.Range("D2:G2").Style = "Normal"

# H2 (rank) is a genuine number.
$q1.Range("H2").Value = 2

# ------------------------------------------------------------------
# 3) Insert a new leading row into "总计" for 2022-Q1 and renumber the
#    index column for the rows that shift down.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy the index cell's style from the row below (still carries the
# original "index column" formatting) onto the newly inserted A2.
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 5.92
$total.Range("B2:D2").Style = "Normal"

# Renumber the rows that were pushed down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
